$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.028692680627711
$ws.Range("C2").Value = 0.03873025157783737
$ws.Range("D2").Value = 0.01263319583351971
$ws.Range("E2").Value = 0.0358191081401645
$ws.Range("F2").Value = 0.05526694318194186
$ws.Range("G2").Value = 0.05181147572192901
$ws.Range("H2").Value = 0.04812787303689326
$ws.Range("I2").Value = 0.03159878169708157
$ws.Range("J2").Value = 0.01155715291314707
$ws.Range("K2").Value = 0.01861302552872287
$ws.Range("L2").Value = 0.038113883352762
$ws.Range("B3").Value = 0.0254213528407077
$ws.Range("C3").Value = 0.02732587332997253
$ws.Range("D3").Value = 0.0202914975688418
$ws.Range("E3").Value = 0.0264414658497892
$ws.Range("F3").Value = 0.0217403166827178
$ws.Range("G3").Value = 0.02464925514519373
$ws.Range("H3").Value = 0.02643469413039877
$ws.Range("I3").Value = 0.0223336180570061
$ws.Range("J3").Value = 0.0078325302161393
$ws.Range("K3").Value = 0.0159358812572107
$ws.Range("L3").Value = 0.0249390324576177
$ws.Range("B4").Value = 0.0215818661329497
$ws.Range("C4").Value = 0.02007100159419284
$ws.Range("D4").Value = 0.0166351930526708
$ws.Range("E4").Value = 0.02098971138387524
$ws.Range("F4").Value = 0.01854822079958457
$ws.Range("G4").Value = 0.02004964564012024
$ws.Range("H4").Value = 0.02118492761448723
$ws.Range("I4").Value = 0.01391382112455054
$ws.Range("J4").Value = 0.005884234056310167
$ws.Range("K4").Value = 0.00158334257103795
$ws.Range("L4").Value = 0.02019421166472205
$ws.Range("B5").Value = 0.0188540837508389
$ws.Range("C5").Value = 0.01742281676510567
$ws.Range("E5").Value = 0.008737176003543963
$ws.Range("F5").Value = 0.01148957549789413
$ws.Range("G5").Value = 0.01598894832380124
$ws.Range("H5").Value = 0.0115275175684763
$ws.Range("I5").Value = 0.0118959649768428
$ws.Range("J5").Value = 0.002741759786351817
$ws.Range("K5").Value = 0.00165935852028116
$ws.Range("L5").Value = 0.01921779057578885
